$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new record row at row 217, pushing the existing rows 217-228
# (and their data) down to 218-229.
$ws.Rows.Item(217).Insert()

# Populate the newly inserted row with the new weekly price record.
$ws.Range("A217").Value = 9
$ws.Range("B217").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C217").Value = "Metropolitana"
$ws.Range("D217").Value = 44585
$ws.Range("E217").Value = 13
$ws.Range("F217").Value = 100112021
$ws.Range("G217").Value = "Ají"
$ws.Range("H217").Value = "Americana (o)"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 16
$ws.Range("K217").Value = 32000
$ws.Range("L217").Value = 34000
$ws.Range("M217").Value = 33000
$ws.Range("N217").Value = "`$/saco 25 kilos"
$ws.Range("O217").Value = "Provincia de Huasco"
$ws.Range("P217").Value = 1320
$ws.Range("Q217").Value = 25
$ws.Range("R217").Value = "Hortaliza"
